# Table2/RRT.xlsx — minor change:
#   - CKD Stage "2.0" sub-row (7 (1.3) / 21 (0.8)) is removed entirely, which
#     shifts every subsequent row up by one (old row 45 disappears, dimension
#     shrinks from D45 to D44).
#   - CKD Stage header row (row 25) values are refreshed: "0.0" -> "0",
#     "469 (90.0)" -> "476 (91.4)", "2411 (93.0)" -> "2432 (93.8)".
#   - What used to be the "3.0" sub-row (old row 27) becomes the new row 26,
#     with its label refreshed from "3.0" -> "3" (counts unchanged).
#   - The two affected merged ranges follow the row shift:
#       A25:A27 -> A25:A26   and   A28:A30 -> A27:A29

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting the entire row 26 shifts rows 27-45 up by one, fixes the
# <dimension> ref automatically, and re-bases the two mergeCell ranges that
# span across the deleted row (A25:A27 -> A25:A26, A28:A30 -> A27:A29).
$ws.Rows.Item(26).Delete()

# Row 25 (CKD Stage header line) gets new figures.
# Column B holds short numeric-looking *text* labels throughout this sheet
# (e.g. "Absent", "2008 - 2010", "3.0", ...), so force text storage here too
# instead of letting COM coerce "0"/"3" into real numbers.
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "0"
$ws.Range("C25").Value = "476 (91.4)"
$ws.Range("D25").Value = "2432 (93.8)"

# New row 26 (was row 27 before the delete) only needs its label tidied up;
# the counts ("45 (8.6)" / "161 (6.2)") are already correct post-shift.
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "3"
